# fix: break line issue
# The salary-input number fields were centre-aligned, which made the
# printed numbers visually "break" against the form's right-aligned
# labels/borders. Right-align all the numeric entry cells, and left-align
# the single "Remarks" (備註) box so wrapped notes read naturally instead
# of being centred.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlLeft  = -4131
$xlRight = -4152

# Numeric / amount entry cells -> right aligned
$rightRanges = @(
    "B7:C7",
    "E7", "E8", "E9", "E10", "E11", "E12", "E13", "E14", "E15", "E16",
    "B17:C17", "E17",
    "B20:C20", "B21:C21", "B22:C22",
    "B24:C24", "B25:C25"
)
foreach ($addr in $rightRanges) {
    $ws.Range($addr).HorizontalAlignment = $xlRight
}

# Remarks box (merged D20:E24) -> left aligned so wrapped text reads naturally
$ws.Range("D20:E24").HorizontalAlignment = $xlLeft

# Leave the cursor where the user finished editing
$ws.Range("E11").Select()
